$d = $word.ActiveDocument

$replacements = @(
    @{old = "401÷6="; new = "703÷5="},
    @{old = "913÷6="; new = "219÷3="},
    @{old = "409÷6="; new = "885÷5="},
    @{old = "198÷4="; new = "109÷8="},
    @{old = "939÷3="; new = "686÷3="},
    @{old = "548÷5="; new = "941÷2="},
    @{old = "274÷4="; new = "637÷7="},
    @{old = "460÷9="; new = "122÷6="},
    @{old = "664÷6="; new = "808÷7="},
    @{old = "852÷2="; new = "900÷6="},
    @{old = "483÷7="; new = "762÷9="},
    @{old = "440÷5="; new = "499÷4="},
    @{old = "561÷2="; new = "501÷7="},
    @{old = "437÷9="; new = "309÷4="},
    @{old = "636÷4="; new = "361÷8="},
    @{old = "145÷9="; new = "847÷3="},
    @{old = "382÷3="; new = "860÷9="},
    @{old = "578÷7="; new = "398÷7="},
    @{old = "255÷3="; new = "742÷4="},
    @{old = "915÷6="; new = "869÷2="},
    @{old = "209÷3="; new = "190÷6="},
    @{old = "304÷7="; new = "228÷8="},
    @{old = "185÷5="; new = "734÷2="},
    @{old = "667÷9="; new = "591÷9="},
    @{old = "806÷2="; new = "239÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
